# 243-peaks37-dno-nomods-lca.xlsx
# - add a new "all fungi" worksheet containing the header row plus every
#   row from the main sheet whose "kingdom" column (D) is "Fungi"
# - add a new, empty "all cyanos" worksheet
# - turn on AutoFilter for the main sheet's data range
# - make "all fungi" the active sheet / tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the two new sheets, right after the main sheet -----------------
$wsFungi = $wb.Worksheets.Add($null, $ws1)
$wsFungi.Name = "all fungi"

$wsCyanos = $wb.Worksheets.Add($null, $wsFungi)
$wsCyanos.Name = "all cyanos"

# --- populate "all fungi" ----------------------------------------------------
# header row
$ws1.Rows(1).Copy()
$wsFungi.Rows(1).PasteSpecial()

# rows on the main sheet where column D ("kingdom") equals "Fungi"
$fungiRows = @(77, 78, 79, 105, 111, 136, 165, 297, 344, 384, 443, 526, 545)
$destRow = 2
foreach ($r in $fungiRows) {
    $ws1.Rows($r).Copy()
    $wsFungi.Rows($destRow).PasteSpecial()
    $destRow++
}

# --- turn on filtering for the main sheet's table ---------------------------
$ws1.Range("A1:AH556").AutoFilter()
$filterName = $ws1.Names.Add("_xlnm._FilterDatabase", "='243-peaks37-dno-nomods-lca'!`$A`$1:`$AH`$556")
$filterName.Visible = $false

# --- page setup tweaks seen in the saved file --------------------------------
$ws1.PageSetup.Orientation = 1
$wsFungi.PageSetup.Orientation = 1

# --- selections / active sheet ----------------------------------------------
$ws1.Range("H5").Select()
$wsFungi.Activate()
$wsFungi.Range("D18").Select()
